# Update "想去人数" (interest count) values in F column across the
# 展览 / 演出 / 全部类型 sheets, matching the refreshed scrape data.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 208
$ws1.Range("F4").Value = 13234
$ws1.Range("F5").Value = 54428
$ws1.Range("F7").Value = 1318
$ws1.Range("F8").Value = 344
$ws1.Range("F10").Value = 866
$ws1.Range("F12").Value = 371
$ws1.Range("F13").Value = 3001
$ws1.Range("F14").Value = 881
$ws1.Range("F15").Value = 5183
$ws1.Range("F16").Value = 1264
$ws1.Range("F17").Value = 967
$ws1.Range("F19").Value = 465
$ws1.Range("F21").Value = 385
$ws1.Range("F22").Value = 1235
$ws1.Range("F24").Value = 34
$ws1.Range("F25").Value = 165
$ws1.Range("F26").Value = 343
$ws1.Range("F27").Value = 7
$ws1.Range("F30").Value = 64
$ws1.Range("F32").Value = 4820
$ws1.Range("F34").Value = 4814
$ws1.Range("F35").Value = 8790
$ws1.Range("F38").Value = 128
$ws1.Range("F40").Value = 415
$ws1.Range("F41").Value = 105
$ws1.Range("F42").Value = 78
$ws1.Range("F43").Value = 4178
$ws1.Range("F44").Value = 210

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 129

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 208
$ws4.Range("F6").Value = 13234
$ws4.Range("F7").Value = 1318
$ws4.Range("F8").Value = 344
$ws4.Range("F10").Value = 866
$ws4.Range("F12").Value = 371
$ws4.Range("F13").Value = 3001
$ws4.Range("F14").Value = 881
$ws4.Range("F16").Value = 1264
$ws4.Range("F18").Value = 129
$ws4.Range("F19").Value = 967
$ws4.Range("F21").Value = 466
$ws4.Range("F22").Value = 385
$ws4.Range("F24").Value = 1235
$ws4.Range("F26").Value = 165
$ws4.Range("F28").Value = 343
$ws4.Range("F31").Value = 4820
$ws4.Range("F33").Value = 4814
$ws4.Range("F34").Value = 8790
$ws4.Range("F37").Value = 128
$ws4.Range("F39").Value = 415
$ws4.Range("F42").Value = 105
$ws4.Range("F43").Value = 78
$ws4.Range("F44").Value = 4178
$ws4.Range("F47").Value = 210
